$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update raw input values; dependent formulas recalc automatically.
$ws.Range("C6").Value = 1350
$ws.Range("D6").Value = 1550
$ws.Range("D8").Value = 1250

# Update the selected cell to match the saved view state.
$ws.Range("E17").Select()
